# "Changes of 27Th April 2022"
#
# The worksheet's C column (ShipmentTrackNum) and the matching D column
# (PackageTrackNum, only present on a handful of rows) get a fresh batch of
# FedEx tracking numbers. The new values must land in the shared-string
# table as plain text (t="s"), exactly like the numbers that were already
# there - not as genuine numbers (Excel would otherwise collapse a value
# like 320018407199 into a Number cell / scientific notation).
#
# Directly assigning Range.Value with a numeric-looking string lets Excel's
# "smart" type inference turn it into a Number. Going through a text
# formula ("=""...""") and then collapsing it to a static value with
# PasteSpecial(xlPasteValues) keeps the literal text but stores it as a
# plain shared string (no stray NumberFormat/quote-prefix style picked up
# along the way).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)   # xlPasteValues
}

# New ShipmentTrackNum for every data row (C2:C22).
$rowValues = [ordered]@{
    2  = "320018407199"
    3  = "320018407203"
    4  = "320018407236"
    5  = "320018407269"
    6  = "320018407306"
    7  = "320018407328"
    8  = "320018407361"
    9  = "320018407394"
    10 = "320018407420"
    11 = "320018407442"
    12 = "320018407486"
    13 = "320018407501"
    14 = "320018407740"
    15 = "320018407773"
    16 = "320018407810"
    17 = "320018407832"
    18 = "320018407876"
    19 = "320018407898"
    20 = "320018407924"
    21 = "320018407946"
    22 = "320018407979"
}

# Rows where column D (PackageTrackNum) mirrors the new column C value.
$dRows = @(5, 6, 7, 13, 14, 15, 16, 17)

foreach ($row in $rowValues.Keys) {
    $val = $rowValues[$row]
    Set-TextValue ($ws.Range("C$row")) $val
    if ($dRows -contains $row) {
        Set-TextValue ($ws.Range("D$row")) $val
    }
}

$excel.CutCopyMode = $false
$ws.Range("C12").Select()
